$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5787152051925659
$ws.Range("B1").Value = 1.18523895740509
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.737319946289062
$ws.Range("E1").Value = 1.454972505569458
